$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.523.61'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '1.881.62'
$ws.Range("E3").Value = '  +1.57%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '0.7148'
$ws.Range("E5").Value = '  +2.29%  '
$ws.Range("D6").Value = '241.93'
$ws.Range("E6").Value = '  +1.75%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").Value = '0.07952'
$ws.Range("E8").Value = '  +1.27%  '
$ws.Range("D9").Value = '0.3113'
$ws.Range("E9").Value = '  +3.13%  '
$ws.Range("D11").Value = '0.08277'
$ws.Range("E11").Value = '  +2.05%  '
$ws.Range("D12").Value = '0.7291'
$ws.Range("E12").Value = '  +3.52%  '
$ws.Range("D13").Value = '5.285'
$ws.Range("E13").Value = '  +2.17%  '
$ws.Range("D14").Value = '1.866.90'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("D15").Value = '91.28'
$ws.Range("E15").Value = '  +2.00%  '
$ws.Range("D16").Value = '29.528.19'
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("D17").Value = '5.945'
$ws.Range("E17").Value = '  +2.48%  '
$ws.Range("D18").Value = '246.46'
$ws.Range("E18").Value = '  +4.67%  '
$ws.Range("D19").Value = '0.000007889'
$ws.Range("E19").Value = '  +0.87%  '
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").Value = '2.130.92'
$ws.Range("E21").Value = '  +2.14%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("D23").Value = '7.972'
$ws.Range("E23").Value = '  +6.22%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '0.1617'
$ws.Range("E25").Value = '  +13.66%  '
$ws.Range("D26").Value = '163.16'
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").Value = '9.073'
$ws.Range("E27").Value = '  +2.34%  '
$ws.Range("D28").Value = '18.35'
$ws.Range("E28").Value = '  +1.91%  '
$ws.Range("D29").Value = '1.359'
$ws.Range("E29").Value = '  -3.38%  '
$ws.Range("D30").Value = '1.502'
$ws.Range("E30").Value = '  +1.98%  '
$ws.Range("D31").Value = '4.393'
$ws.Range("E31").Value = '  +1.60%  '
$ws.Range("D32").Value = '4.105'
$ws.Range("E32").Value = '  +2.43%  '
$ws.Range("D33").Value = '0.05279'
$ws.Range("E33").Value = '  +2.44%  '
$ws.Range("D34").Value = '1.963'
$ws.Range("E34").Value = '  +2.33%  '
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("D36").Value = '0.7279'
$ws.Range("E36").Value = '  +2.51%  '
$ws.Range("D37").Value = '2.679'
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '0.01871'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = '1.232.13'
$ws.Range("E39").Value = '  +7.07%  '
$ws.Range("D40").Value = '2.717'
$ws.Range("E40").Value = '  +0.39%  '
$ws.Range("D41").Value = '0.9132'
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("D42").Value = '6.216'
$ws.Range("E42").Value = '  +4.30%  '
$ws.Range("D43").Value = '73.94'
$ws.Range("E43").Value = '  +5.50%  '
$ws.Range("E44").Value = '  +0.11%  '
$ws.Range("D45").Value = '102.31'
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("D46").Value = '2.026.61'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").Value = '0.5290'
$ws.Range("E47").Value = '  -0.09%  '
$ws.Range("D48").Value = '1.801'
$ws.Range("E48").Value = '  +3.78%  '
$ws.Range("D49").Value = '2.942'
$ws.Range("E49").Value = '  +11.35%  '
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").Value = '9.325'
$ws.Range("E51").Value = '  +2.30%  '
